$d = $word.ActiveDocument

# --- Change 1: "Definice done" -> "Definice " + "DONE" (two runs, same formatting) ---
$titlePara = $d.Paragraphs(1).Range
$titleText = $titlePara.Text
$idx = $titleText.IndexOf("done")
if ($idx -ge 0) {
    $start = $titlePara.Start + $idx
    $end = $start + 4

    # Replace the lowercase "done" with uppercase "DONE" text.
    $wordRange = $d.Range($start, $end)
    $wordRange.Text = "DONE"

    # Force Word to split this into its own run (distinct from the
    # preceding "Definice " run) by toggling a character property and
    # then reverting it back to the original value. The net formatting
    # is unchanged, but the run boundary introduced by the edit remains.
    $doneRange = $d.Range($start, $start + 4)
    $doneRange.Bold = $true
    $doneRange2 = $d.Range($start, $start + 4)
    $doneRange2.Bold = $false
}

# --- Change 2: merge "Prošel i" + "nterní Review" runs into one run ---
$d.Content.Find.Execute("nterní Review", $true, $false, $false, $false, $false, $true, 1, $false, "nterní Review", 2) | Out-Null
